# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values on the active worksheet for rows 3-12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    3  = 5
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 2
    9  = 1
    10 = 3
    11 = 0
    12 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
